# Update crypto price/volume snapshot values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text is purely numeric-looking ("325.56", "5.530", "0.00001033", ...)
# must be forced to Text format first, otherwise Excel auto-converts the assigned
# string into a real number (losing trailing zeros / switching to scientific notation)
# same as the original inline-string cells, which must stay text.
$textPriceCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D11",
    "D13",
    "D14",
    "D15",
    "D18",
    "D19",
    "D20",
    "D22",
    "D26",
    "D28",
    "D29",
    "D30",
    "D32",
    "D33",
    "D34",
    "D35",
    "D37",
    "D38",
    "D39",
    "D40",
    "D42",
    "D43",
    "D44",
    "D45",
    "D47",
    "D48",
    "D49",
    "D51"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.482.53"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.905.21"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "325.56"
$ws.Range("E5").Value = "  -2.63%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "0.4801"
$ws.Range("E7").Value = "  +2.48%  "
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").Value = "0.08075"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").Value = "23.42"
$ws.Range("E11").Value = "  +4.55%  "
$ws.Range("D12").Value = "1.912.42"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "5.958"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").Value = "7.085"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").Value = "90.16"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "0.00001033"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "17.65"
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "29.485.33"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "5.547"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").Value = "2.127.97"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("D26").Value = "154.61"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "6.115"
$ws.Range("E28").Value = "  +6.01%  "
$ws.Range("D29").Value = "2.099"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("D30").Value = "118.42"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("E31").Value = "  -3.32%  "
$ws.Range("D32").Value = "0.09516"
$ws.Range("D33").Value = "5.530"
$ws.Range("E33").Value = "  +2.03%  "
$ws.Range("D34").Value = "1.396"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("D35").Value = "3.544"
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").Value = "0.06080"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").Value = "1.178"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "0.5904"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "7.937"
$ws.Range("E40").Value = "  -5.71%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").Value = "10.23"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "1.283"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.07807"
$ws.Range("E44").Value = "  +3.94%  "
$ws.Range("D45").Value = "2.401"
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").Value = "0.5535"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("D48").Value = "1.926"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").Value = "114.13"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").Value = "72.34"
$ws.Range("E51").Value = "  +0.90%  "
